# Fix Training Data Issue
# Data was taken from 1 day off due to way NBA stats were shown.
# The "Date" column (BF) held the literal text "5-25-2012-13" on every
# data row; it should read "2013-05-25" instead (ISO formatted date,
# stored as plain text - NOT an Excel date serial number).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldText = "5-25-2012-13"
$newText = "2013-05-25"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

# Locate the "Date" column dynamically from the header row (row 1)
$dateCol = $null
$headerFound = $ws.Rows.Item(1).Find("Date")
if ($headerFound -ne $null) {
    $dateCol = $headerFound.Column
} else {
    $dateCol = 58 # fallback: column BF
}

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $dateCol)
    $current = $cell.Value()
    if ($current -eq $oldText) {
        # Force a text number format first so Excel does not reinterpret
        # the replacement string "2013-05-25" as a date serial value.
        $cell.NumberFormat = "@"
        $cell.Value = $newText
    }
}
